$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 565.2222
$ws.Range("I58").Value = 385.875
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1157.625
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -1007.625
$ws.Range("N58").Value = -6300
$ws.Range("H88").Value = 2908.5557
$ws.Range("I88").Value = 3167.6667
$ws.Range("J88").Value = 2779
$ws.Range("K88").Value = 3167.6667
$ws.Range("L88").Value = 2779
$ws.Range("M88").Value = -2761.6667
$ws.Range("N88").Value = -3591
$ws.Range("H91").Value = 2908.5557
$ws.Range("I91").Value = 3167.6667
$ws.Range("J91").Value = 2779
$ws.Range("K91").Value = 3167.6667
$ws.Range("L91").Value = 2779
$ws.Range("M91").Value = -1763.6667
$ws.Range("N91").Value = -5587

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2367.4736
$ws.Range("I2").Value = 1824.5
$ws.Range("K2").Value = 1824.5
$ws.Range("M2").Value = -1711.5
$ws.Range("H74").Value = 2676.5454
$ws.Range("I74").Value = 2629.2
$ws.Range("K74").Value = 2629.2
$ws.Range("M74").Value = -1755.2
$ws.Range("H76").Value = 32000
$ws.Range("J76").Value = 32000
$ws.Range("L76").Value = 32000
$ws.Range("N76").Value = -32676
$ws.Range("H77").Value = 2676.5454
$ws.Range("I77").Value = 2629.2
$ws.Range("K77").Value = 13146
$ws.Range("M77").Value = -8778
$ws.Range("H79").Value = 32000
$ws.Range("J79").Value = 32000
$ws.Range("L79").Value = 32000
$ws.Range("N79").Value = -34340
$ws.Range("H116").Value = 2367.4736
$ws.Range("I116").Value = 1824.5
$ws.Range("K116").Value = 1824.5
$ws.Range("M116").Value = 469.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2367.4736
$ws.Range("I3").Value = 1824.5
$ws.Range("K3").Value = 1824.5
$ws.Range("M3").Value = -1710.5
$ws.Range("H86").Value = 3505.5
$ws.Range("I86").Value = 3505.5
$ws.Range("K86").Value = 3505.5
$ws.Range("M86").Value = -2382.5
$ws.Range("H89").Value = 3505.5
$ws.Range("I89").Value = 3505.5
$ws.Range("K89").Value = 17527.5
$ws.Range("M89").Value = -11911.5
$ws.Range("H134").Value = 2896.6956
$ws.Range("I134").Value = 2137.3684
$ws.Range("K134").Value = 6412.1052
$ws.Range("M134").Value = -3877.1052

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("K16").Value = 5000
$ws.Range("M16").Value = -4713
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 350
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H81").Value = 2840.4546
$ws.Range("J81").Value = 3498.75
$ws.Range("L81").Value = 10496.25
$ws.Range("N81").Value = -12742.25
$ws.Range("H84").Value = 2840.4546
$ws.Range("J84").Value = 3498.75
$ws.Range("L84").Value = 31488.75
$ws.Range("N84").Value = -42720.75
$ws.Range("H121").Value = 10000938
$ws.Range("I121").Value = 520
$ws.Range("J121").Value = 15001148
$ws.Range("K121").Value = 1560
$ws.Range("L121").Value = 45003444
$ws.Range("M121").Value = -250
$ws.Range("N121").Value = -45006064
$ws.Range("H132").Value = 557948.9
$ws.Range("I132").Value = 2188.75
$ws.Range("J132").Value = 1669469.1
$ws.Range("K132").Value = 19698.75
$ws.Range("L132").Value = 15025221.9
$ws.Range("M132").Value = -17168.75
$ws.Range("N132").Value = -15030281.9

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 8043.6
$ws.Range("I2").Value = 27.75
$ws.Range("K2").Value = 27.75
$ws.Range("M2").Value = 85.25
$ws.Range("H70").Value = 22932.969
$ws.Range("J70").Value = 4480.769
$ws.Range("L70").Value = 4480.769
$ws.Range("N70").Value = -5020.769
$ws.Range("H73").Value = 22932.969
$ws.Range("J73").Value = 4480.769
$ws.Range("L73").Value = 4480.769
$ws.Range("N73").Value = -6352.769
$ws.Range("H102").Value = 3287.5
$ws.Range("I102").Value = 3050
$ws.Range("K102").Value = 3050
$ws.Range("M102").Value = -1428
$ws.Range("H132").Value = 2950.0833
$ws.Range("J132").Value = 5338
$ws.Range("L132").Value = 16014
$ws.Range("N132").Value = -21074

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2003
$ws.Range("I100").Value = 2003
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2003
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1462
$ws.Range("N100").ClearContents()
$ws.Range("H123").Value = 29888
$ws.Range("J123").Value = 29888
$ws.Range("L123").Value = 29888
$ws.Range("N123").Value = -39688
$ws.Range("H132").Value = 39379.9
$ws.Range("I132").Value = 48045.348
$ws.Range("J132").Value = 6162.3335
$ws.Range("K132").Value = 144136.044
$ws.Range("L132").Value = 18487.0005
$ws.Range("M132").Value = -141606.044
$ws.Range("N132").Value = -23547.0005

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3851.3333
$ws.Range("I81").Value = 2874.5557
$ws.Range("J81").Value = 5316.5
$ws.Range("K81").Value = 5749.1114
$ws.Range("L81").Value = 10633
$ws.Range("M81").Value = -4688.1114
$ws.Range("N81").Value = -12755
$ws.Range("H84").Value = 3851.3333
$ws.Range("I84").Value = 2874.5557
$ws.Range("J84").Value = 5316.5
$ws.Range("K84").Value = 28745.557
$ws.Range("L84").Value = 53165
$ws.Range("M84").Value = -23441.557
$ws.Range("N84").Value = -63773
$ws.Range("H113").Value = 583.53845
$ws.Range("J113").Value = 690
$ws.Range("L113").Value = 2070
$ws.Range("N113").Value = -6410
$ws.Range("H115").Value = 57051
$ws.Range("J115").Value = 57051
$ws.Range("L115").Value = 57051
$ws.Range("N115").Value = -60185
$ws.Range("H132").Value = 2986.8125
$ws.Range("I132").Value = 2922.3845
$ws.Range("J132").Value = 3266
$ws.Range("K132").Value = 8767.1535
$ws.Range("L132").Value = 9798
$ws.Range("M132").Value = -6237.1535
$ws.Range("N132").Value = -14858
$ws.Range("H136").Value = 13335599
$ws.Range("I136").Value = 1798.9
$ws.Range("K136").Value = 5396.700000000001
$ws.Range("M136").Value = -2846.700000000001
